$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(3)
$shape.TextFrame.TextRange.Paragraphs(1).Font.Size = 20
